# Fix plot dimensions on page 2 (Sheet1):
#  - update the image-filename strings in column D to the current
#    auto-generated plot filenames
#  - update the "w" (width) column (E) values used to size the plots
#  - move the active selection to E6 to match where editing left off

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Commercial_LONGFINSQUID_Landings_LBS_2025-04-07.png"
$ws.Range("D3").Value = "N_Commercial_Vessels_Landing_LONGFINSQUID_2025-04-07.png"
$ws.Range("D4").Value = "TOTALANNUALREV_LONGFINSQUID_2023Dols_2025-04-07.png"
$ws.Range("D5").Value = "western gulf stream index_2025-04-07.png"

$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3

[void]$ws.Range("E6").Select()
